$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark that currently sits after "Unique Poudel".
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2. Split the sentence "... will help data scientists with the deep
#    understanding ... solution." into two runs and re-insert the "_GoBack"
#    bookmark (as a zero-length / collapsed bookmark) between them, right
#    before "deep understanding".
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("data scientists with the ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)

# Re-write the text of the now-isolated trailing run in place so the engine
# recomputes its xml:space bookkeeping cleanly (it otherwise inherits the
# leading run's "preserve" flag verbatim).
$r2 = $d.Content
$tail = "deep understanding of the craft of the problem and plan formulation to engineer the solution."
$null = $r2.Find.Execute($tail, $true, $false, $false, $false, $false, $true, 1, $false, $tail, 2)

# ---------------------------------------------------------------------------
# 3. Turn the single list item "Help documents, content-creation on Surround,
#    website design, UX/UI requirements (if any)" into two list items:
#      - "Logo Design which reflects the Surround"
#      - "Help documents, content-creation on Surround, website design,
#         UX/UI requirements "
# ---------------------------------------------------------------------------
$old3 = "Help documents, content-creation on Surround, website design, UX/UI requirements (if any)"
$new3 = "Logo Design which reflects the Surround" + [char]13 + "Help documents, content-creation on Surround, website design, UX/UI requirements "
$r3 = $d.Content
$null = $r3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# The paragraph-break insert above does not copy the run-level language
# formatting onto the new second paragraph's run, so restore it explicitly.
$r4 = $d.Content
$tail3 = "Help documents, content-creation on Surround, website design, UX/UI requirements "
$null = $r4.Find.Execute($tail3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4.LanguageID = "en-AU"
$r4.LanguageIDFarEast = "en-AU"
